$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: property_category column (I) should read "building", not "land"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# 汽車 (car) sheet: property_category column (H) should read "car", not "land"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
